$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.832.19"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.707.73"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'318.07"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "'0.3966"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("D9").Value = "'1.510"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").Value = "'1.006"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").Value = "'53.25"
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("D12").Value = "'0.08963"
$ws.Range("E12").Value = "  +2.31%  "
$ws.Range("D13").Value = "'7.719"
$ws.Range("E13").Value = "  +6.83%  "
$ws.Range("D14").Value = "'24.41"
$ws.Range("E14").Value = "  +5.05%  "
$ws.Range("D15").Value = "'8.183"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("D17").Value = "1.716.01"
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").Value = "'100.27"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'0.07165"
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("D20").Value = "'20.09"
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("D21").Value = "'7.515"
$ws.Range("E21").Value = "  +6.76%  "
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").Value = "'14.51"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("D24").Value = "24.838.81"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").Value = "'3.105"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").Value = "'2.342"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").Value = "'23.06"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("D28").Value = "'9.273"
$ws.Range("E28").Value = "  +23.36%  "
$ws.Range("D29").Value = "'166.43"
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("D30").Value = "'139.78"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("D31").Value = "'5.227"
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("D32").Value = "'7.906"
$ws.Range("E32").Value = "  +10.73%  "
$ws.Range("D33").Value = "'0.09099"
$ws.Range("E33").Value = "  +6.21%  "
$ws.Range("D34").Value = "'1.085"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").Value = "'0.03078"
$ws.Range("E35").Value = "  +12.89%  "
$ws.Range("E36").Value = "  +3.12%  "
$ws.Range("D37").Value = "'11.15"
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("D38").Value = "'1.970"
$ws.Range("E38").Value = "  +2.17%  "
$ws.Range("D39").Value = "'14.61"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("D40").Value = "'0.09315"
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("D41").Value = "'1.488"
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").Value = "'0.7865"
$ws.Range("E42").Value = "  +3.07%  "
$ws.Range("D43").Value = "'16.54"
$ws.Range("E43").Value = "  +5.78%  "
$ws.Range("D44").Value = "'2.648"
$ws.Range("E44").Value = "  +3.45%  "
$ws.Range("D45").Value = "'0.7324"
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("D46").Value = "'4.258"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").Value = "'1.008"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").Value = "'1.354"
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("D49").Value = "'141.11"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "'94.96"
$ws.Range("E50").Value = "  +5.72%  "
$ws.Range("D51").Value = "'0.08070"
$ws.Range("E51").Value = "  +1.13%  "
